$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '33.069.69'
$ws.Range("E2").Value = '  +10.60%  '
$ws.Range("D3").Value = '1.737.67'
$ws.Range("E3").Value = '  +4.99%  '
$ws.Range("D4").Value = '''0.987'
$ws.Range("E4").Value = '  -0.99%  '
$ws.Range("D5").Value = '''225.05'
$ws.Range("E5").Value = '  +3.64%  '
$ws.Range("D6").Value = '''0.541'
$ws.Range("E6").Value = '  +4.00%  '
$ws.Range("D7").Value = '''0.989'
$ws.Range("E7").Value = '  -0.81%  '
$ws.Range("D8").Value = '''31.85'
$ws.Range("E8").Value = '  +9.36%  '
$ws.Range("D9").Value = '''45.40'
$ws.Range("E9").Value = '  +3.81%  '
$ws.Range("D10").Value = '''0.275'
$ws.Range("E10").Value = '  +4.86%  '
$ws.Range("D11").Value = '''0.0662'
$ws.Range("E11").Value = '  +8.31%  '
$ws.Range("D12").Value = '''0.0913'
$ws.Range("E12").Value = '  +1.50%  '
$ws.Range("D13").Value = '1.973.28'
$ws.Range("E13").Value = '  +4.68%  '
$ws.Range("D14").Value = '1.745.72'
$ws.Range("E14").Value = '  +5.50%  '
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").Value = '''0.629'
$ws.Range("E15").Value = '  +5.29%  '
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").Value = '''10.51'
$ws.Range("E16").Value = '  +5.36%  '
$ws.Range("D17").Value = '32.960.53'
$ws.Range("E17").Value = '  +10.17%  '
$ws.Range("D18").Value = '''4.24'
$ws.Range("E18").Value = '  +7.81%  '
$ws.Range("D19").Value = '''68.28'
$ws.Range("E19").Value = '  +5.70%  '
$ws.Range("D20").Value = '''256.55'
$ws.Range("E20").Value = '  +6.80%  '
$ws.Range("D21").Value = '0.0₃0735'
$ws.Range("E21").Value = '  +3.85%  '
$ws.Range("D22").Value = '''0.979'
$ws.Range("E22").Value = '  -1.77%  '
$ws.Range("D23").Value = '''10.40'
$ws.Range("E23").Value = '  +4.11%  '
$ws.Range("D24").Value = '''4.33'
$ws.Range("E24").Value = '  +4.20%  '
$ws.Range("D25").Value = '''2.15'
$ws.Range("E25").Value = '  -1.08%  '
$ws.Range("D26").Value = '''158.58'
$ws.Range("E26").Value = '  +0.49%  '
$ws.Range("D27").Value = '''16.45'
$ws.Range("E27").Value = '  +4.66%  '
$ws.Range("E28").Value = '  +3.58%  '
$ws.Range("D29").Value = '''6.91'
$ws.Range("E29").Value = '  +2.77%  '
$ws.Range("D30").Value = '''0.988'
$ws.Range("E30").Value = '  -0.89%  '
$ws.Range("D31").Value = '''3.82'
$ws.Range("E31").Value = '  +12.68%  '
$ws.Range("D32").Value = '''0.0512'
$ws.Range("E32").Value = '  +2.90%  '
$ws.Range("D33").Value = '''1.17'
$ws.Range("E33").Value = '  +5.35%  '
$ws.Range("D34").Value = '''3.45'
$ws.Range("E34").Value = '  +7.32%  '
$ws.Range("D35").Value = '1.549.65'
$ws.Range("E35").Value = '  +8.19%  '
$ws.Range("D36").Value = '''1.78'
$ws.Range("E36").Value = '  +4.85%  '
$ws.Range("D37").Value = '''84.92'
$ws.Range("E37").Value = '  +9.59%  '
$ws.Range("D38").Value = '''1.04'
$ws.Range("E38").Value = '  +1.49%  '
$ws.Range("D39").Value = '''0.625'
$ws.Range("E39").Value = '  +9.07%  '
$ws.Range("E40").Value = '  +5.24%  '
$ws.Range("D41").Value = '''2.69'
$ws.Range("E41").Value = '  +0.80%  '
$ws.Range("E42").Value = '  +0.12%  '
$ws.Range("D43").Value = '''2.09'
$ws.Range("E43").Value = '  +7.16%  '
$ws.Range("E44").Value = '  +2.77%  '
$ws.Range("D45").Value = '''0.0503'
$ws.Range("E45").Value = '  +0.30%  '
$ws.Range("D46").Value = '''55.29'
$ws.Range("E46").Value = '  +9.57%  '
$ws.Range("E47").Value = '  +4.42%  '
$ws.Range("D48").Value = '1.877.05'
$ws.Range("E48").Value = '  +4.78%  '
$ws.Range("D49").Value = '''0.987'
$ws.Range("E49").Value = '  -0.96%  '
$ws.Range("D50").Value = '''5.61'
$ws.Range("E50").Value = '  +4.69%  '
$ws.Range("D51").Value = '''94.79'
$ws.Range("E51").Value = '  +0.66%  '
